# Regenerate merged AHB files
#
# 1) Rename the "_old" / "_new" header-suffix columns to "_FV2404" / "_FV2410"
# 2) Turn the A1:U79 range into an Excel Table ("Table1")
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) held the "*_old" headers -> becomes "*_FV2404"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2404"
}

# Column K (11) is "diff" and is left untouched.

# Columns L..U (12..21) held the "*_new" headers -> becomes "*_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2410"
}

# Convert the data range into a table, matching the existing used range.
$tableRange = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Applied FV2410/FV2404 header rename, Table1, and frozen header row."
